$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A4: "Non existent" -> "Invalid"
$ws.Range("A4").Value = "Invalid"

# Select A4 (matches the selection change seen in the diff)
$ws.Range("A4").Select()
